$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 84, pushing existing rows 84-199 down to 85-200.
$ws.Rows("84:84").Insert()

# Populate the newly inserted row 84 with the new weekly data point.
$ws.Range("A84").Value = 3
$ws.Range("B84").Value = "Femacal de La Calera"
$ws.Range("C84").Value = "Coquimbo"
$ws.Range("D84").Value = 44482
$ws.Range("D84").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E84").Value = 5
$ws.Range("F84").Value = 100112012
$ws.Range("G84").Value = "Espinaca"
$ws.Range("H84").Value = "Sin especificar"
$ws.Range("I84").Value = "Primera"
$ws.Range("J84").Value = 160
$ws.Range("K84").Value = 3500
$ws.Range("L84").Value = 3500
$ws.Range("M84").Value = 3500
$ws.Range("N84").Value = "$/docena de atados (3 kilos)"
$ws.Range("O84").Value = "Provincia de Quillota"
$ws.Range("P84").Value = 1167
$ws.Range("Q84").Value = 3
$ws.Range("R84").Value = "Hortaliza"
